# "List of parts updated"
#
# - Row 42: STATUS changes from "Ordered" to "Ready"
# - New rows 43-51: nine new parts (fuse holder, ARK connector, fuse,
#   and a set of SMD resistors), each with qty/status/price/date/link/seller
# - "Main" table grows from A1:H42 to A1:H51
# - J2 (PROJECT PRICE) recalculates automatically via its SUM formula

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42: mark the ZS-X10 driver as received ---
$ws.Range("C42").Value = "Ready"

# --- New parts data (row, name, qty, status, unit price, date serial, link, seller) ---
$rows = @(
    @{ R = 43; Name = "PTF-76 fuse holder";                      Qty = 1;  Price = 1.85; Date = 45267; Link = "https://allegro.pl/oferta/oprawa-bezpiecznika-5x20mm-ptf-76-z-uchwytem-1727-10001770448"; Seller = "Allegro (tomsyty)" }
    @{ R = 44; Name = "ARK 2 PIN connector";                     Qty = 30; Price = 0.6;  Date = 45268; Link = "https://allegro.pl/oferta/zlacze-ark-2-pin-5-08mm-5-szt-0693-5-10912515733"; Seller = "Allegro (tomsyty)" }
    @{ R = 45; Name = "3A fuse (5 x 20mm, fast, 10 pcs.)";       Qty = 1;  Price = 1.85; Date = 45269; Link = "https://allegro.pl/oferta/bezpiecznik-5x20mm-szybki-3a-10-szt-0602-10-9510464377"; Seller = "Allegro (tomsyty)" }
    @{ R = 46; Name = "SMD resistor (1206, 100 Ohm, 50 pcs.)";   Qty = 1;  Price = 2;    Date = 45270; Link = "https://allegro.pl/oferta/rezystory-smd-1206-100r-101-50-szt-4425-50-12684715501"; Seller = "Allegro (tomsyty)" }
    @{ R = 47; Name = "SMD resistor (1206, 680 Ohm, 50 pcs.)";   Qty = 1;  Price = 2;    Date = 45271; Link = "https://allegro.pl/oferta/rezystory-smd-1206-680r-681-50-szt-4439-50-12684756050"; Seller = "Allegro (tomsyty)" }
    @{ R = 48; Name = "SMD resistor (1206, 1000 Ohm, 50 pcs.)";  Qty = 1;  Price = 2;    Date = 45272; Link = "https://allegro.pl/oferta/rezystory-smd-1206-1k-102-50-szt-4441-50-12684760810"; Seller = "Allegro (tomsyty)" }
    @{ R = 49; Name = "SMD resistor (1206, 1500 Ohm, 50 pcs.)";  Qty = 1;  Price = 2;    Date = 45273; Link = "https://allegro.pl/oferta/rezystory-smd-1206-1-5k-152-50-szt-4443-50-12684768996"; Seller = "Allegro (tomsyty)" }
    @{ R = 50; Name = "SMD resistor (1206, 3300 Ohm, 50 pcs.)";  Qty = 1;  Price = 2;    Date = 45274; Link = "https://allegro.pl/oferta/rezystory-smd-1206-3-3k-332-50-szt-4449-50-12684791144"; Seller = "Allegro (tomsyty)" }
    @{ R = 51; Name = "SMD resistor (1206, 4700 Ohm, 50 pcs.)";  Qty = 1;  Price = 2.05; Date = 45275; Link = "https://allegro.pl/oferta/rezystory-smd-1206-4-7k-472-50-szt-4452-50-12684801225"; Seller = "Allegro (tomsyty)" }
)

foreach ($row in $rows) {
    $r = $row.R

    $ws.Range("A$r").Value = $row.Name
    $ws.Range("B$r").Value = $row.Qty
    $ws.Range("C$r").Value = "Ordered"
    $ws.Range("D$r").Value = $row.Price
    $ws.Range("E$r").Formula = "=PRODUCT(B$r*D$r)"
    $ws.Range("F$r").Value = $row.Date

    $ws.Range("G$r").Value = $row.Link
    $ws.Hyperlinks.Add($ws.Range("G$r"), $row.Link)
    $ws.Range("G$r").Style = "Hyperlink"

    $ws.Range("H$r").Value = $row.Seller
}

# --- Grow the "Main" table to cover the new rows ---
$lo = $ws.ListObjects.Item("Main")
$lo.Resize($ws.Range("A1:H51"))

# --- Move the active selection to reflect the new bottom of the sheet ---
$ws.Range("H52").Select()
